$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 15:11"

$ws.Range("B4").Value = 7947549
$ws.Range("C4").Value = 2044
$ws.Range("D4").Value = 5090226
$ws.Range("E4").Value = 2638027
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 219296

$ws.Range("B22").Value = 339267
$ws.Range("C22").Value = 323
$ws.Range("D22").Value = 325330
$ws.Range("E22").Value = 8894
$ws.Range("G22").Value = 25
$ws.Range("H22").Value = 5043

$ws.Range("B30").Value = 174653
$ws.Range("C30").Value = 6373
$ws.Range("G30").Value = 17
$ws.Range("H30").Value = 6584

$ws.Range("B41").Value = 111116
$ws.Range("C41").Value = 548
$ws.Range("D41").Value = 103268
$ws.Range("E41").Value = 7190
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 658

$ws.Range("B51").Value = 86664
$ws.Range("C51").Value = 1090
$ws.Range("D51").Value = 53187
$ws.Range("E51").Value = 31397
$ws.Range("G51").Value = 13
$ws.Range("H51").Value = 2080

$ws.Range("B78").Value = 32422
$ws.Range("C78").Value = 340
$ws.Range("D78").Value = 26380
$ws.Range("E78").Value = 5373
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 669

$ws.Range("B107").Value = 10222
$ws.Range("C107").Value = 42
$ws.Range("D107").Value = 9058
$ws.Range("E107").Value = 1085

$ws.Range("B175").Value = 554
$ws.Range("C175").Value = 5
$ws.Range("E175").Value = 15
